# "Default values for excel"
# Replace the example file-path / text placeholders on the Meteo_Input and
# Soil_Input sheets with concrete numeric default values, and update the
# sheet selections / active tab to reflect the state the workbook was left
# in when it was saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Meteo_Input: row 2 held example string paths in several columns; these
# are replaced with representative numeric defaults.
# ---------------------------------------------------------------------
$wsMeteo = $wb.Worksheets.Item("Meteo_Input")

$wsMeteo.Range("B2").Value = 25.8
$wsMeteo.Range("C2").Value = 20.93
$wsMeteo.Range("D2").Value = 44
$wsMeteo.Range("E2").Value = 70.22
$wsMeteo.Range("G2").Value = 0.3
$wsMeteo.Range("H2").Value = 0.38
$wsMeteo.Range("J2").Value = 310
$wsMeteo.Range("M2").Value = 569

# ---------------------------------------------------------------------
# Soil_Input: same idea - row 2 string placeholders become numeric
# defaults, and the trailing unused cells (J2:O2) are cleared out
# completely.
# ---------------------------------------------------------------------
$wsSoil = $wb.Worksheets.Item("Soil_Input")

$wsSoil.Range("B2").Value = 0.4
$wsSoil.Range("C2").Value = 0.4
$wsSoil.Range("F2").Value = 0.32
$wsSoil.Range("G2").Value = 0.04
$wsSoil.Range("J2:O2").Clear()

# ---------------------------------------------------------------------
# Update the selections left on each sheet, and which sheet/tab is
# active when the workbook is reopened.
# ---------------------------------------------------------------------
[void]$wsSoil.Range("C8").Select()

$wsAdditional = $wb.Worksheets.Item("Additional_Input")
[void]$wsAdditional.Range("F2").Select()

[void]$wsMeteo.Activate()
[void]$wsMeteo.Range("D8").Select()
